$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$labels = @(
    "Danh mục lương",
    "Tổng công tại CẦN THƠ",
    "Lương cơ bản tại CẦN THƠ",
    "Chiết khấu sale chính tại CẦN THƠ",
    "Chiết khấu sale phụ tại CẦN THƠ",
    "Đơn 1 bác sĩ tại CẦN THƠ",
    "Đơn 2 bác sĩ tại CẦN THƠ",
    "Công phụ phẫu 1 tại CẦN THƠ",
    "Công phụ phẫu 2 tại CẦN THƠ",
    "Ứng lương tại CẦN THƠ",
    "Tổng công tại LONG XUYÊN",
    "Lương công tác tại LONG XUYÊN",
    "Lương cơ bản tại LONG XUYÊN",
    "Chiết khấu sale chính tại LONG XUYÊN",
    "Chiết khấu sale phụ tại LONG XUYÊN",
    "Đơn 1 bác sĩ tại LONG XUYÊN",
    "Đơn 2 bác sĩ tại LONG XUYÊN",
    "Công phụ phẫu 1 tại LONG XUYÊN",
    "Công phụ phẫu 2 tại LONG XUYÊN",
    "Ứng lương tại LONG XUYÊN",
    "Tổng công tại SÓC TRĂNG",
    "Lương công tác tại SÓC TRĂNG",
    "Lương cơ bản tại SÓC TRĂNG",
    "Chiết khấu sale chính tại SÓC TRĂNG",
    "Chiết khấu sale phụ tại SÓC TRĂNG",
    "Đơn 1 bác sĩ tại SÓC TRĂNG",
    "Đơn 2 bác sĩ tại SÓC TRĂNG",
    "Công phụ phẫu 1 tại SÓC TRĂNG",
    "Công phụ phẫu 2 tại SÓC TRĂNG",
    "Ứng lương tại SÓC TRĂNG",
    "Tổng lương tại CẦN THƠ",
    "Tổng lương tại LONG XUYÊN",
    "Tổng lương tại SÓC TRĂNG",
    "Tổng lương"
)

$values = @(
    1,
    28,
    2060000,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    3090000,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    3090000,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    2060000,
    3090000,
    3090000,
    8239999.999999998
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
